$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# The sheet ships protected; temporarily unprotect so the cell values below
# can be written, then restore protection afterwards.
$ws.Unprotect()

# Update the confidential disclaimer date from 2021-05-24 to 2021-05-25
$ws.Range("A9").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-05-25 for illustrative purposes only and are subject to change."
# Re-fit the row height so it matches the sheet's default (setting a
# multi-line value otherwise leaves an explicit custom row height behind).
$ws.Rows.Item(9).AutoFit()

# Update the Weight (D) and Percent Change (E) values for rows 2-6
$ws.Range("D2").Value = 0.2551014941107076
$ws.Range("E2").Value = -0.008469726450726855

$ws.Range("D3").Value = 0.2544469428103732
$ws.Range("E3").Value = -0.01006355932203384

$ws.Range("D4").Value = 0.2426963351702937
$ws.Range("E4").Value = -0.002989969135802517

$ws.Range("D5").Value = 0.2477552279086254
$ws.Range("E5").Value = -0.0002555583950932805

$ws.Range("E6").Value = -0.005510252255523707

# Restore sheet protection
$ws.Protect()
